$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the original authoring order so shared-string indices line up:
# E1, E2, F1, F2, then fill down E3:E21, then fill down F3:F21.
$ws.Range("E1").Value = "Unit ID"
$ws.Range("E2").Value = "asdf123"
$ws.Range("F1").Value = "Owner"
$ws.Range("F2").Value = "Nice"

for ($row = 3; $row -le 21; $row++) {
    $unitNum = 123 + ($row - 2)
    $ws.Cells.Item($row, 5).Value = "asdf$unitNum"
}

for ($row = 3; $row -le 21; $row++) {
    if ($row -le 13) {
        $ws.Cells.Item($row, 6).Value = "Nice"
    } else {
        $ws.Cells.Item($row, 6).Value = "Not Nice "
    }
}

# Widen column D to fit the newly highlighted content
$ws.Columns.Item(4).ColumnWidth = 55.75

# Update the active selection to reflect where the user ended up
$ws.Range("E12").Select()
